# Update gh-pages to output generated at 456a3b4
#
# The upstream scraper re-ran and produced a slightly newer snapshot:
#  - one cancelled event ("苏州·首届 童年回忆同人only 茶歇聚会（取消）", 2024-10-05)
#    dropped out of the "展览" (exhibition) listing - and therefore also out
#    of the combined "全部类型" (all types) listing - so every later row
#    shifts up by one.
#  - the "想去人数" (interest-count) numbers ticked up for most events.
#  - the "演出" (performance) sheet's one live count also ticked up.
#
# Sheet layout (tab order): 展览(1), 演出(2), 本地生活(3), 全部类型(4)
#
# NOTE: this runtime's PowerShell function calls only bind parameters
# positionally (named "-param value" binding does not work), so
# Set-Row below is always called with every positional slot, using
# $null for "leave this column alone".

$wb = $excel.ActiveWorkbook

function Set-Row {
    param($ws, $row, $b, $c, $d, $e, $f, $g, $h, $i)
    $ws.Cells.Item($row, 1).Value = $row - 1
    if ($null -ne $b) { $ws.Cells.Item($row, 2).Value = $b }
    if ($null -ne $c) { $ws.Cells.Item($row, 3).Value = $c }
    if ($null -ne $d) { $ws.Cells.Item($row, 4).Value = $d }
    if ($null -ne $e) { $ws.Cells.Item($row, 5).Value = $e }
    if ($null -ne $f) { $ws.Cells.Item($row, 6).Value = $f }
    if ($null -ne $g) { $ws.Cells.Item($row, 7).Value = $g }
    if ($null -ne $h) { $ws.Cells.Item($row, 8).Value = $h }
    if ($null -ne $i) { $ws.Cells.Item($row, 9).Value = $i }
}

# ---------------------------------------------------------------------------
# Sheet 1: 展览 (Exhibitions)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Drop the cancelled "童年回忆" teahouse meetup row; everything below moves up.
$ws1.Rows.Item(3).Delete()

# Refresh the "想去人数" (F) / "最低票价" (G) counters that changed, and make
# sure the leading index column (A) stays a plain 1..N sequence after the
# shift (row 2 -> A=1, row 3 -> A=2, ...).
Set-Row $ws1 2  $null $null $null $null 780   $null          $null $null
Set-Row $ws1 3  $null $null $null $null 55    59             $null $null
Set-Row $ws1 4  $null $null $null $null 411   "不可售"        $null $null
Set-Row $ws1 5  $null $null $null $null 136   139            $null $null
Set-Row $ws1 6  $null $null $null $null 15    125            $null $null
Set-Row $ws1 7  $null $null $null $null 155   70             $null $null
Set-Row $ws1 8  $null $null $null $null 340   58             $null $null
Set-Row $ws1 9  $null $null $null $null 452   78             $null $null
Set-Row $ws1 10 $null $null $null $null 510   178            $null $null
Set-Row $ws1 11 $null $null $null $null 139   40             $null $null
Set-Row $ws1 12 $null $null $null $null 11718 60             $null $null
Set-Row $ws1 13 $null $null $null $null 5406  55             $null $null

# ---------------------------------------------------------------------------
# Sheet 2: 演出 (Performances) - no rows added/removed, just a counter bump.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
Set-Row $ws2 2 $null $null $null $null 105 $null $null $null

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活 (Local life) - untouched (header only).
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (All types) - same cancelled-row drop as 展览, plus the
# same counter refresh across the combined listing (展览 + 演出 rows).
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Rows.Item(3).Delete()

Set-Row $ws4 2  $null $null $null $null 780   $null          $null $null
Set-Row $ws4 3  $null $null $null $null 55    59             $null $null
Set-Row $ws4 4  $null $null $null $null 105   88             $null $null
Set-Row $ws4 5  $null $null $null $null 2     220            $null $null
Set-Row $ws4 6  $null $null $null $null 411   "不可售"        $null $null
Set-Row $ws4 7  $null $null $null $null 136   139            $null $null
Set-Row $ws4 8  $null $null $null $null 15    125            $null $null
Set-Row $ws4 9  $null $null $null $null 155   70             $null $null
Set-Row $ws4 10 $null $null $null $null 340   58             $null $null
Set-Row $ws4 11 $null $null $null $null 452   78             $null $null
Set-Row $ws4 12 $null $null $null $null 510   178            $null $null
Set-Row $ws4 13 $null $null $null $null 139   40             $null $null
Set-Row $ws4 14 $null $null $null $null 11718 60             $null $null
Set-Row $ws4 15 $null $null $null $null 7     280            $null $null
Set-Row $ws4 16 $null $null $null $null 5406  55             $null $null
